# REST API set up
# Row 5 used to be employee "Wick, John" (3165 / 1ABC126); the record was
# changed to employee 53331 "Smith, John" (keeping the same license plate).
# Removing the last reference to "Wick, John" drops that now-unused shared
# string automatically when the workbook is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("A5").Value = 53331
$ws.Range("B5").Value = "Smith, John"

# Move the active selection to B5 (and clear the prior scrolled-to-G1 view)
# to mirror the author's final cursor position when they saved the file.
[void]$ws.Range("B5").Select()
